# repull data, push all data, mean calculation
# Update the dSF column (F) values for a set of rows to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    11 = 2
    13 = -1
    15 = 1
    18 = -2
    19 = 9
    20 = -1
    24 = -6
    28 = 2
    29 = -2
    36 = 4
    37 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
